$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from the Chinese default name to "Sheet1"
$ws.Name = "Sheet1"

# Move the active selection to C8 (was H2)
$ws.Range("C8").Select()
